# Lab1/Statistcs.xlsx - "Fix a bug. Make rat avoid walls. Update a graphical performance"
#
# The maze simulation's wall/mud density results changed: the "Miss" column is
# now always 0/0/0 (the rat no longer crashes into walls), and the "Moves" /
# "Stucks" density triples were re-measured. Update the data table, the
# selection/scroll view and the new columns' widths to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data: re-run statistics (Miss is now always 0/0/0 - rat avoids walls) ---
# Write the "Miss" column (E) first, then walk rows 7 -> 2 for Moves (D) /
# Stucks (F) so new shared-string entries land in the same order as the
# target workbook.
$ws.Range("E2").Value = "0/0/0"
$ws.Range("E3").Value = "0/0/0"
$ws.Range("E4").Value = "0/0/0"
$ws.Range("E5").Value = "0/0/0"
$ws.Range("E6").Value = "0/0/0"
$ws.Range("E7").Value = "0/0/0"

$ws.Range("D7").Value = "58/802/60"
$ws.Range("F7").Value = "38/233/22"

$ws.Range("D6").Value = "184/100/14"
$ws.Range("F6").Value = "99/114/4"

$ws.Range("D5").Value = "94/318/66"
$ws.Range("F5").Value = "31/246/12"

$ws.Range("D4").Value = "32/40/46"
$ws.Range("F4").Value = "12/33/21"

$ws.Range("D3").Value = "24/10/12"
$ws.Range("F3").Value = "14/0/0"

$ws.Range("D2").Value = "2/4/4"
$ws.Range("F2").Value = "0/0/2"

# --- Graphical performance / view update ---
# New columns D:F got their own (wider) custom widths instead of inheriting
# the sheet default.
$ws.Columns.Item(4).ColumnWidth = 13.833333333333332
$ws.Columns.Item(5).ColumnWidth = 15.499999999999998
$ws.Columns.Item(6).ColumnWidth = 15.499999999999998

# Scroll the sheet over and move the selection from G4 (off the data table)
# to F4, matching the new, narrower data range.
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("F4").Select() | Out-Null
